$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 85
$ws.Range("F6").Value = 326
$ws.Range("F10").Value = 1333
$ws.Range("F13").Value = 174
$ws.Range("F16").Value = 121
$ws.Range("F17").Value = 259
$ws.Range("F18").Value = 1685
$ws.Range("F21").Value = 255
$ws.Range("F22").Value = 2870
$ws.Range("F23").Value = 25
$ws.Range("F24").Value = 407
$ws.Range("F29").Value = 2851
$ws.Range("F30").Value = 1648
$ws.Range("F33").Value = 683
$ws.Range("F34").Value = 870
$ws.Range("F35").Value = 1903
$ws.Range("F37").Value = 1912
$ws.Range("F39").Value = 36
$ws.Range("F40").Value = 49
$ws.Range("F42").Value = 50
$ws.Range("F43").Value = 900
$ws.Range("F44").Value = 811
$ws.Range("F45").Value = 1049
$ws.Range("F46").Value = 128
$ws.Range("F47").Value = 449
$ws.Range("F48").Value = 229
$ws.Range("F49").Value = 3368

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 810

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 85
$ws.Range("F7").Value = 326
$ws.Range("F11").Value = 1333
$ws.Range("F14").Value = 174
$ws.Range("F17").Value = 121
$ws.Range("F18").Value = 259
$ws.Range("F19").Value = 1685
$ws.Range("F22").Value = 255
$ws.Range("F23").Value = 2870
$ws.Range("F24").Value = 25
$ws.Range("F25").Value = 407
$ws.Range("F28").Value = 2851
$ws.Range("F29").Value = 1648
$ws.Range("F33").Value = 810
$ws.Range("F35").Value = 870
$ws.Range("F36").Value = 1903
$ws.Range("F39").Value = 1912
$ws.Range("F41").Value = 900
$ws.Range("F42").Value = 811
$ws.Range("F43").Value = 1049
$ws.Range("F44").Value = 128
$ws.Range("F45").Value = 449
$ws.Range("F47").Value = 229
$ws.Range("F48").Value = 3368

